$wb = $excel.ActiveWorkbook

# Rename the first sheet from "$100" to "$125 " (note trailing space)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "$125 "

# Work on the "$125" sheet (second sheet) - add a new delivery row
$ws2 = $wb.Worksheets.Item(2)

# Copy formatting from the last existing row (118) down into the new row (119)
$ws2.Range("A118:C118").Copy()
$ws2.Range("A119:C119").PasteSpecial(-4122)

$ws2.Cells.Item(119, 1).Value = 20637
$ws2.Cells.Item(119, 2).Value = "Hughesville"
$ws2.Cells.Item(119, 3).Value = "MD"
$ws2.Rows.Item(119).RowHeight = 15.75

# Update the remembered selection/scroll position on the "$125" sheet
$null = $ws2.Range("D114").Select()

# Restore the original active sheet so the workbook still opens on "$125 "
$null = $ws1.Activate()
